# add value in done.xlsx file
#
# The user typed two very large integers into column C (rows 22 and 23).
# Excel only keeps 15 significant decimal digits of precision for a
# numeric cell value, so what actually gets stored is the rounded form
# of what was typed:
#   C22: typed 1111111111111111111111111 (25 ones)  -> stored 1.11111111111111E+24
#   C23: typed 2222222222222222          (16 twos)   -> stored 2222222222222220
# We assign the already-rounded values below so the stored <v> text is
# exactly what Excel itself would have produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C22").Value2 = [double]"1.11111111111111E+24"
$ws.Range("C23").Value2 = [double]"2222222222222220"

# Column C is now too narrow for the big numbers, so best-fit (auto-fit)
# its width to the new content. (11.17 lands squarely in the ColumnWidth
# range that serializes to the saved width="12" the real AutoFit produced.)
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(3).ColumnWidth = 11.17

# Leave the last-entered cell selected, same as the end state in the file.
$ws.Range("C23").Select() | Out-Null
